$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$loginValid = "Login with valid username and password"
$loginCred = "Login using credentials from Excel with Apache POI"
$passed = "PASSED"
$edge = "edge"

$ws.Range("A2").Value = $loginValid
$ws.Range("B2").Value = $passed
$ws.Range("C2").Value = $edge

$ws.Range("A3").Value = $loginCred
$ws.Range("B3").Value = $passed
$ws.Range("C3").Value = $edge

$ws.Range("A4").Value = $loginValid
$ws.Range("B4").Value = $passed
$ws.Range("C4").Value = $edge

$ws.Range("A5").Value = $loginCred
$ws.Range("B5").Value = $passed
$ws.Range("C5").Value = $edge
